$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 222.85715
$ws.Range("J6").Value = 174.5
$ws.Range("L6").Value = 523.5
$ws.Range("N6").Value = -747.5

# Row 19
$ws.Range("H19").Value = 1788.8
$ws.Range("I19").Value = 1765
$ws.Range("K19").Value = 1765
$ws.Range("M19").Value = -1590

# Row 51
$ws.Range("H51").Value = 9219
$ws.Range("I51").Value = 9044.75
$ws.Range("K51").Value = 9044.75
$ws.Range("M51").Value = -8560.75

# Row 86
$ws.Range("H86").Value = 3294.0417
$ws.Range("I86").Value = 3040.4443
$ws.Range("J86").Value = 3446.2
$ws.Range("K86").Value = 3040.4443
$ws.Range("L86").Value = 3446.2
$ws.Range("M86").Value = -1917.4443
$ws.Range("N86").Value = -5692.2

# Row 89
$ws.Range("H89").Value = 3294.0417
$ws.Range("I89").Value = 3040.4443
$ws.Range("J89").Value = 3446.2
$ws.Range("K89").Value = 15202.2215
$ws.Range("L89").Value = 17231
$ws.Range("M89").Value = -9586.2215
$ws.Range("N89").Value = -28463

# Row 92
$ws.Range("H92").Value = 3346
$ws.Range("I92").Value = 1296.25
$ws.Range("K92").Value = 1296.25
$ws.Range("M92").Value = -48.25

# Row 98
$ws.Range("H98").Value = 1640.625
$ws.Range("I98").Value = 1160.7142
$ws.Range("K98").Value = 1160.7142
$ws.Range("M98").Value = 337.2858000000001

# Row 101
$ws.Range("H101").Value = 10180.929
$ws.Range("I101").Value = 8359.429
$ws.Range("J101").Value = 12002.429
$ws.Range("K101").Value = 25078.287
$ws.Range("L101").Value = 36007.287
$ws.Range("M101").Value = -23456.287
$ws.Range("N101").Value = -39251.287

# Row 109
$ws.Range("H109").Value = 93496
$ws.Range("J109").Value = 93496
$ws.Range("L109").Value = 93496
$ws.Range("N109").Value = -96270

# Row 116
$ws.Range("H116").Value = 6990.6665
$ws.Range("I116").Value = 6990.6665
$ws.Range("K116").Value = 6990.6665
$ws.Range("M116").Value = -3548.6665

# Row 122
$ws.Range("H122").Value = 1640.625
$ws.Range("I122").Value = 1160.7142
$ws.Range("K122").Value = 3482.1426
$ws.Range("M122").Value = -1032.1426

# Row 132
$ws.Range("H132").Value = 914072.0600000001
$ws.Range("I132").Value = 5099.125
$ws.Range("K132").Value = 15297.375
$ws.Range("M132").Value = -12767.375

# Row 137
$ws.Range("H137").Value = 1212.1464
$ws.Range("I137").Value = 852.05554
$ws.Range("K137").Value = 2556.16662
$ws.Range("M137").Value = -6.166619999999966

# Row 141
$ws.Range("H141").Value = 5837.385
$ws.Range("I141").Value = 3642.3333
$ws.Range("K141").Value = 10926.9999
$ws.Range("M141").Value = -5746.999899999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4271.0757
$ws.Range("I32").Value = 1624.8298
$ws.Range("K32").Value = 1624.8298
$ws.Range("M32").Value = -1337.8298

# Row 102
$ws.Range("H102").Value = 1958
$ws.Range("I102").Value = 1958
$ws.Range("K102").Value = 1958
$ws.Range("M102").Value = -336

# Row 114
$ws.Range("H114").Value = 22000
$ws.Range("J114").Value = 22000
$ws.Range("L114").Value = 22000
$ws.Range("N114").Value = -30678

$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 16250
$ws.Range("I11").Value = 16250
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 16250
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -16110
$ws.Range("N11").Value = $null

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2134.6155
$ws.Range("I31").Value = 1535.5
$ws.Range("J31").Value = 2833.5833
$ws.Range("K31").Value = 1535.5
$ws.Range("L31").Value = 2833.5833
$ws.Range("M31").Value = -1240.5
$ws.Range("N31").Value = -3423.5833

# Row 34
$ws.Range("H34").Value = 2134.6155
$ws.Range("I34").Value = 1535.5
$ws.Range("J34").Value = 2833.5833
$ws.Range("K34").Value = 1535.5
$ws.Range("L34").Value = 2833.5833
$ws.Range("M34").Value = -1333.5
$ws.Range("N34").Value = -3237.5833

# Row 86
$ws.Range("H86").Value = 4999.6665
$ws.Range("I86").Value = 4499.5
$ws.Range("K86").Value = 4499.5
$ws.Range("M86").Value = -3376.5

# Row 89
$ws.Range("H89").Value = 4999.6665
$ws.Range("I89").Value = 4499.5
$ws.Range("K89").Value = 22497.5
$ws.Range("M89").Value = -16881.5

# Row 105
$ws.Range("H105").Value = 2679
$ws.Range("I105").Value = 1518.5
$ws.Range("K105").Value = 1518.5
$ws.Range("M105").Value = 228.5

# Row 107
$ws.Range("H107").Value = 1457.6
$ws.Range("I107").Value = 1167.6154
$ws.Range("K107").Value = 1167.6154
$ws.Range("M107").Value = 752.3846000000001

# Row 134
$ws.Range("H134").Value = 3150
$ws.Range("I134").Value = 3150
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9450
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6915
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 288.95
$ws.Range("I12").Value = 390.55554
$ws.Range("J12").Value = 205.81818
$ws.Range("K12").Value = 1171.66662
$ws.Range("L12").Value = 617.4545400000001
$ws.Range("M12").Value = -998.66662
$ws.Range("N12").Value = -963.4545400000001

# Row 50
$ws.Range("H50").Value = 3458.1365
$ws.Range("I50").Value = 849.5
$ws.Range("J50").Value = 3719
$ws.Range("K50").Value = 2548.5
$ws.Range("L50").Value = 11157
$ws.Range("M50").Value = -2067.5
$ws.Range("N50").Value = -12119

# Row 53
$ws.Range("H53").Value = 3458.1365
$ws.Range("I53").Value = 849.5
$ws.Range("J53").Value = 3719
$ws.Range("K53").Value = 2548.5
$ws.Range("L53").Value = 11157
$ws.Range("M53").Value = -2067.5
$ws.Range("N53").Value = -12119

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 1156.2
$ws.Range("I107").Value = 450.66666
$ws.Range("K107").Value = 450.66666
$ws.Range("M107").Value = 1469.33334

# Row 113
$ws.Range("H113").Value = 3061
$ws.Range("I113").Value = 2730.5
$ws.Range("K113").Value = 2730.5
$ws.Range("M113").Value = -560.5

# Row 126
$ws.Range("H126").Value = 4698.8335
$ws.Range("J126").Value = 4698.8335
$ws.Range("L126").Value = 14096.5005
$ws.Range("N126").Value = -19036.5005

# Row 132
$ws.Range("H132").Value = 5486.1665
$ws.Range("I132").Value = 6213.3335
$ws.Range("J132").Value = 4759
$ws.Range("K132").Value = 18640.0005
$ws.Range("L132").Value = 14277
$ws.Range("M132").Value = -16110.0005
$ws.Range("N132").Value = -19337

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 325.55554
$ws.Range("I16").Value = 325.55554
$ws.Range("K16").Value = 325.55554
$ws.Range("M16").Value = -155.55554

# Row 61
$ws.Range("H61").Value = 1630.5
$ws.Range("I61").Value = 1630.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1630.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1428.5
$ws.Range("N61").Value = $null

# Row 113
$ws.Range("H113").Value = 1630.5
$ws.Range("I113").Value = 1630.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1630.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 539.5
$ws.Range("N113").Value = $null

# Row 122
$ws.Range("H122").Value = 3494.6956
$ws.Range("I122").Value = 3596.7693
$ws.Range("J122").Value = 3362
$ws.Range("K122").Value = 10790.3079
$ws.Range("L122").Value = 10086
$ws.Range("M122").Value = -8340.3079
$ws.Range("N122").Value = -14986

# Row 132
$ws.Range("H132").Value = 3625.7144
$ws.Range("I132").Value = 3127.6667
$ws.Range("K132").Value = 9383.000100000001
$ws.Range("M132").Value = -6853.000100000001

# Row 136
$ws.Range("H136").Value = 1506
$ws.Range("J136").Value = 3165.3333
$ws.Range("L136").Value = 9495.999899999999
$ws.Range("N136").Value = -14595.9999

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 2428024.2
$ws.Range("I96").Value = 3033530.2
$ws.Range("J96").Value = 6000
$ws.Range("K96").Value = 3033530.2
$ws.Range("L96").Value = 6000
$ws.Range("M96").Value = -3032157.2
$ws.Range("N96").Value = -8746

# Row 113
$ws.Range("H113").Value = 426.4

# Row 122
$ws.Range("H122").Value = 4910.878
$ws.Range("I122").Value = 5378.7036
$ws.Range("K122").Value = 16136.1108
$ws.Range("M122").Value = -13686.1108

# Row 132
$ws.Range("H132").Value = 9000
$ws.Range("I132").Value = 9000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 27000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -24470
$ws.Range("N132").Value = $null
